$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values stored as literal text (matches source data,
# which sometimes uses thousand-separator dots, e.g. "42.418.55"). Force the
# cell number format to Text before assignment so Excel does not silently
# reinterpret values like "0.630" or "0.0000105" as numbers (which would
# drop significant trailing/representational digits).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.418.55'
$ws.Range("E2").Value = '  -1.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.335.97'
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.59'
$ws.Range("E5").Value = '  -2.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.41'
$ws.Range("E6").Value = '  -0.89%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("E7").Value = '  -2.04%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.611'
$ws.Range("E9").Value = '  -8.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.73'
$ws.Range("E10").Value = '  -2.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0917'
$ws.Range("E11").Value = '  -2.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.34'
$ws.Range("E12").Value = '  -2.89%  '
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.983'
$ws.Range("E14").Value = '  -4.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.77'
$ws.Range("E15").Value = '  -6.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.693.67'
$ws.Range("E16").Value = '  -1.71%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.364.78'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.418.43'
$ws.Range("E18").Value = '  -1.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.54'
$ws.Range("E19").Value = '  -4.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000105'
$ws.Range("E20").Value = '  -2.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '75.55'
$ws.Range("E21").Value = '  -1.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.56'
$ws.Range("E22").Value = '  +5.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '254.17'
$ws.Range("E23").Value = '  -7.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.28'
$ws.Range("E24").Value = '  -5.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.28'
$ws.Range("E25").Value = '  -2.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.29'
$ws.Range("E27").Value = '  -3.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.66'
$ws.Range("E28").Value = '  -2.64%  '
$ws.Range("E29").Value = '  +2.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.33'
$ws.Range("E30").Value = '  +0.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.96'
$ws.Range("E31").Value = '  -5.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0881'
$ws.Range("E32").Value = '  -4.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.01'
$ws.Range("E33").Value = '  +2.69%  '
$ws.Range("E34").Value = '  -8.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.123'
$ws.Range("E35").Value = '  +15.64%  '
$ws.Range("E36").Value = '  -2.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.56'
$ws.Range("E37").Value = '  -6.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0359'
$ws.Range("E38").Value = '  -2.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.81'
$ws.Range("E39").Value = '  -9.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.63'
$ws.Range("E40").Value = '  -7.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.234'
$ws.Range("E41").Value = '  +1.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '70.16'
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.45'
$ws.Range("E43").Value = '  -9.09%  '
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.82'
$ws.Range("E45").Value = '  -5.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '111.61'
$ws.Range("E46").Value = '  -9.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.44'
$ws.Range("E47").Value = '  -1.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '84.92'
$ws.Range("E48").Value = '  -10.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.04'
$ws.Range("E49").Value = '  -4.47%  '
$ws.Range("E50").Value = '  -2.91%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0997'
$ws.Range("E51").Value = '  -2.11%  '
